$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.291.84'
$ws.Range("E2").Value = '  +5.10%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.914.07'
$ws.Range("E3").Value = '  +5.54%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9996'
$ws.Range("E4").Value = '  +0.11%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '253.30'
$ws.Range("E5").Value = '  +0.85%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9996'
$ws.Range("E6").Value = '  +0.09%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5149'
$ws.Range("E7").Value = '  +2.88%  '

$ws.Range("E8").Value = '  +6.20%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2984'
$ws.Range("E9").Value = '  +7.40%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06794'
$ws.Range("E10").Value = '  +6.29%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.914.52'
$ws.Range("E11").Value = '  +5.56%  '

$ws.Range("E12").Value = '  +3.75%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07374'
$ws.Range("E13").Value = '  +3.08%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6989'
$ws.Range("E14").Value = '  +7.10%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '87.69'
$ws.Range("E15").Value = '  +6.96%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.900'
$ws.Range("E16").Value = '  +3.92%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.299.90'
$ws.Range("E17").Value = '  +5.23%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008020'
$ws.Range("E18").Value = '  +8.33%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9991'
$ws.Range("E19").Value = '  +0.05%  '

$ws.Range("E20").Value = '  +6.00%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.162.23'
$ws.Range("E21").Value = '  +5.76%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9985'
$ws.Range("E22").Value = '  +0.03%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.851'
$ws.Range("E23").Value = '  +4.92%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.740'
$ws.Range("E24").Value = '  +7.36%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.171'
$ws.Range("E25").Value = '  +3.03%  '

$ws.Range("B26").Value = 'BitcoinCash'
$ws.Range("C26").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '139.41'
$ws.Range("E26").Value = '  +21.87%  '

$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '146.46'
$ws.Range("E27").Value = '  +1.77%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.25'
$ws.Range("E28").Value = '  +7.59%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.016'
$ws.Range("E29").Value = '  +6.76%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.397'
$ws.Range("E30").Value = '  +0.46%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.261'
$ws.Range("E31").Value = '  +2.00%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08829'
$ws.Range("E32").Value = '  +5.50%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.018'
$ws.Range("E33").Value = '  +4.29%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05111'
$ws.Range("E34").Value = '  +2.60%  '

$ws.Range("E35").Value = '  +6.04%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7208'
$ws.Range("E36").Value = '  +6.09%  '

$ws.Range("E37").Value = '  +0.09%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.845'
$ws.Range("E38").Value = '  +3.93%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.320'
$ws.Range("E39").Value = '  +5.99%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9763'
$ws.Range("E40").Value = '  +0.61%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.01696'
$ws.Range("E41").Value = '  +6.40%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.081'
$ws.Range("E42").Value = '  +1.45%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '106.10'
$ws.Range("E43").Value = '  +4.33%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4319'
$ws.Range("E44").Value = '  +4.42%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9989'

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.722'
$ws.Range("E46").Value = '  +6.61%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1284'
$ws.Range("E47").Value = '  +4.61%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05760'
$ws.Range("E48").Value = '  +4.65%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '33.43'
$ws.Range("E49").Value = '  +5.44%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.535'
$ws.Range("E50").Value = '  +3.79%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3833'
$ws.Range("E51").Value = '  +4.84%  '
